$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price history rows (newest first), to be written into A2:C15.
# Row1 (header: Price Date / Price / Currency - Unit Level) is unchanged.
$data = @(
    @("04/11/2025", "0.953", "SGD"),
    @("03/11/2025", "0.955", "SGD"),
    @("31/10/2025", "0.953", "SGD"),
    @("30/10/2025", "0.953", "SGD"),
    @("29/10/2025", "0.954", "SGD"),
    @("28/10/2025", "0.955", "SGD"),
    @("27/10/2025", "0.955", "SGD"),
    @("24/10/2025", "0.953", "SGD"),
    @("23/10/2025", "0.952", "SGD"),
    @("22/10/2025", "0.952", "SGD"),
    @("21/10/2025", "0.951", "SGD"),
    @("17/10/2025", "0.950", "SGD"),
    @("16/10/2025", "0.949", "SGD"),
    @("15/10/2025", "0.949", "SGD")
)

$lastRow = 1 + $data.Count

# Format the destination range as Text first so date-/number-looking
# strings ("04/11/2025", "0.953", ...) are kept as literal text instead
# of being auto-converted into date serials / numbers.
$targetRange = $ws.Range("A2:C$lastRow")
$targetRange.NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $ws.Range("A$row").Value = $data[$i][0]
    $ws.Range("B$row").Value = $data[$i][1]
    $ws.Range("C$row").Value = $data[$i][2]
}

# Reset cell formatting back to the sheet's default (plain, General)
# style now that the literal text is safely stored, so the cells don't
# keep an explicit "@" number format applied to them.
$targetRange.ClearFormats()
